$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Förändrad" (column C) changed from 45205 to 45206 for every existing
#    data row (rows 2..291).
for ($r = 2; $r -le 291; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}

# 2. Row 291 gains an explicit "ht=15 customHeight=1" row-height flag
#    (cosmetic row-height normalisation applied by the source tool).
$ws.Rows.Item(291).RowHeight = 15

# 3. A brand-new row (292) is appended for case "A 47905-2023".
$newRow = 292
$ws.Cells.Item($newRow, 1).Value = "A 47905-2023"

$ws.Cells.Item($newRow, 2).Value = 45204
$ws.Cells.Item($newRow, 2).NumberFormat = $ws.Cells.Item($newRow - 1, 2).NumberFormat

$ws.Cells.Item($newRow, 3).Value = 45206
$ws.Cells.Item($newRow, 3).NumberFormat = $ws.Cells.Item($newRow - 1, 3).NumberFormat

$ws.Cells.Item($newRow, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item($newRow, 5).Value = "ORSA"
$ws.Cells.Item($newRow, 6).Value = "Allmännings- och besparingsskogar"
$ws.Cells.Item($newRow, 7).Value = 9.300000000000001
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0

# Column R ("Artnamn") is always present but wrap-text styled even when
# empty - mirror that by copying the style from the row above.
$ws.Cells.Item($newRow, 18).Value = ""
$ws.Cells.Item($newRow, 18).WrapText = $true

Write-Output "edit applied"
